# Automatische test-sync: 2025-06-30 20:08:50
#
# Adds a new "Testmail #14" log entry to the Logs sheet, updates the
# conditional formatting ranges to cover the new row, swaps the order of
# two existing Dashboard categories, appends a new Dashboard category row
# for "Intern verzoek / Actie voor medewerker", and extends the bar chart's
# category/value references to include that new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 14 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Graag contact opnemen met de klant hierover."
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Testmail #14: Graag contact opnemen met de klant hierover."
$logs.Range("D14").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E14").Value = "Beste klant,`nBedankt voor uw bericht. Kunt u meer informatie geven over waarover u graag contact wilt opnemen? Bijvoorbeeld over welke specifieke kwestie of vraag het gaat? Met deze aanvullende details kunnen we u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F14").Value = "2025-06-30 20:08:20"
$logs.Range("G14").Value = "Ja"
$logs.Range("H14").Value = "Nee"
$logs.Range("I14").Value = "Ja"
$logs.Range("J14").Value = "Nee"

# Setting a multi-line value (via embedded newline) auto-marks the row
# with a custom height; re-running AutoFit restores the implicit/default
# row height so the row matches the style of the other data rows.
$logs.Rows.Item(14).AutoFit()

# Extend the conditional formatting ranges from row 13 to row 14, keeping
# priorities / dxfIds / rules intact.
$dRange = $logs.Range("D2:D13")
for ($i = 1; $i -le $dRange.FormatConditions.Count; $i++) {
    $dRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D14"))
}

$gRange = $logs.Range("G2:G13")
for ($i = 1; $i -le $gRange.FormatConditions.Count; $i++) {
    $gRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G14"))
}

$hRange = $logs.Range("H2:H13")
for ($i = 1; $i -le $hRange.FormatConditions.Count; $i++) {
    $hRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("H2:H14"))
}

$iRange = $logs.Range("I2:I13")
for ($i = 1; $i -le $iRange.FormatConditions.Count; $i++) {
    $iRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("I2:I14"))
}

$jRange = $logs.Range("J2:J13")
for ($i = 1; $i -le $jRange.FormatConditions.Count; $i++) {
    $jRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("J2:J14"))
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: swap rows 4 and 5, append new row 8
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Openingstijden / Locatie"
$dash.Range("A5").Value = "Bestelling / Levering"

$dash.Range("A8").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B8").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series references to row 8
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$8,'Dashboard'!`$B`$2:`$B`$8,1)"
